$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix the wording of the existing "disable functionality" MOM item on the
#    20-04-2022 sheet (shared string text edit).
# ---------------------------------------------------------------------------
$ws20 = $wb.Worksheets.Item("20-04-2022")
$ws20.Range("E17").Value = "Don't use remove functionality , instead use disable functionality"

# ---------------------------------------------------------------------------
# 2) Create the new "21-04-2022" MOM sheet. The easiest way to reproduce the
#    exact same look & feel (borders / fills / merged cells) as the previous
#    day's sheet is to duplicate it and then edit the duplicate's content.
# ---------------------------------------------------------------------------
$ws20.Copy([System.Reflection.Missing]::Value, $ws20)
$ws21 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws21.Name = "21-04-2022"

# The new sheet's content starts two rows lower than 20-04-2022's did, so
# shift everything down by inserting two blank rows above the table.
$ws21.Rows("6:7").Insert()

# ---------------------------------------------------------------------------
# 3) Replace the Q&A content with the 21-04-2022 meeting minutes.
# ---------------------------------------------------------------------------
$ws21.Range("E9").Value = "Need one more Attribute in both Organisation and Department entity (Operation for getting list of departments under one specific organisation)"
$ws21.Range("F9").Value = "Yes"

$ws21.Range("E11").Value = "Multiple  filters simultaneously"
$ws21.Range("F11").Value = "No"

$ws21.Range("E13").Value = "GetByStatus , PostView redundancy"
$ws21.Range("F13").Value = "Yes"

$ws21.Range("E15").Value = "Naming convention - Award Request become Award and Award becomes Award Type"
$ws21.Range("F15").Value = "Yes"

$ws21.Range("E17").Value = "Write operation in Status Entity"
$ws21.Range("F17").Value = "Yes"

$ws21.Range("E19").Value = "Need Designation ID, Department ID, ACE ID in Employee Entity"
$ws21.Range("F19").Value = "Yes"

$ws21.Range("E21").Value = "Logout's  input and output  method"
$ws21.Range("F21").Value = ""

$ws21.Range("E23").Value = "Pass objects instead of IDs"
$ws21.Range("F23").Value = "Yes"

# ---------------------------------------------------------------------------
# 4) Restore/adjust the view state: the new sheet becomes the active tab,
#    20-04-2022 keeps a plain view, and 12-04-2022 scrolls further down.
# ---------------------------------------------------------------------------
$ws21.Range("E15:E16").Select()
$ws21.Activate()

$ws20.Range("E17:E18").Select()

$ws12 = $wb.Worksheets.Item("12-04-2022")
$ws12.Range("A13").Select()

$ws21.Activate()
